$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6, copying row 5's formatting (styles, font size, number format)
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(6).Insert(-4121)  # xlShiftDown
$ws.Range("F6").Clear()

# Fill in the new product row values
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Lavadora 15kg"
$ws.Range("C6").Value = 10000
$ws.Range("D6").Value = "Línea Blanca"
$ws.Range("E6").Value = "Lavadora automática de 15 kg."

# Match row height used by the rest of the table
$ws.Rows.Item(6).RowHeight = 28.5

# Match selection from the diff
$ws.Range("D5").Select()
